$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row (county list) based on column B.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

# Header for new column C
$ws.Range("C1").Value = "desc_upper"

# Fill column C with the uppercase of column B for each data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $b = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 3).Value = $b.ToUpper()
}

# Update the selection to match the authored state.
$ws.Range("B15").Select()
